$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the password for the Noor.Uddin.* rows (B2:B4) to the new value.
$ws.Range("B2").Value = "MHRA12345"
$ws.Range("B3").Value = "MHRA12345"
$ws.Range("B4").Value = "MHRA12345"

# Move the selection to B4 as recorded in the saved view state.
$ws.Activate()
$ws.Range("B4").Select()
